# step_2 updated and reran
# Append a duplicate of the "DiSCoVER: top drugs (cerebellar stem cell
# control)" slide (slide 6) to the end of the deck, as a new slide 9.

$p = $ppt.ActivePresentation

# Slide 6 holds the canonical "DiSCoVER: top drugs" title + table that the
# new trailing slide should reproduce exactly.
$src = $p.Slides.Item(6)

# Duplicate() clones the slide (shapes, table, formatting) verbatim and
# inserts the copy immediately after the source slide.
$new = $src.Duplicate()

# Move the freshly duplicated slide to the very end of the deck so it
# becomes the last slide (new sldId appended after the existing ones).
$new.MoveTo($p.Slides.Count)
